# Auto-generated edit script: updates numeric columns H-N for specific rows
# across all 8 worksheets, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 100
$ws.Range("N4").ClearContents()
$ws.Range("J4").Value = 0
$ws.Range("I4").Value = 100
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 14
$ws.Range("K4").Value = 100
$ws.Range("L40").Value = 52633696
$ws.Range("K40").Value = 1500
$ws.Range("H40").Value = 41668656
$ws.Range("N40").Value = -52634046
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 52633696
$ws.Range("M40").Value = -1325
$ws.Range("L137").Value = 88238538
$ws.Range("J137").Value = 29412846
$ws.Range("N137").Value = -88243638
$ws.Range("M137").Value = -1154.4489
$ws.Range("H137").Value = 7576953
$ws.Range("I137").Value = 1234.8163
$ws.Range("K137").Value = 3704.4489

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -15863.378
$ws.Range("I32").Value = 16150.378
$ws.Range("N32").Value = -68298
$ws.Range("L32").Value = 67724
$ws.Range("J32").Value = 67724
$ws.Range("K32").Value = 16150.378
$ws.Range("H32").Value = 18549.15
$ws.Range("L45").Value = 4752.3335
$ws.Range("I45").Value = 8729.714
$ws.Range("M45").Value = -8352.714
$ws.Range("N45").Value = -5506.3335
$ws.Range("H45").Value = 7536.5
$ws.Range("J45").Value = 4752.3335
$ws.Range("K45").Value = 8729.714
$ws.Range("K110").Value = 654.3333
$ws.Range("I110").Value = 654.3333
$ws.Range("H110").Value = 681.2778
$ws.Range("M110").Value = 1390.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I105").Value = 4101.6665
$ws.Range("L105").Value = 3000
$ws.Range("J105").Value = 3000
$ws.Range("N105").Value = -6494
$ws.Range("M105").Value = -2354.6665
$ws.Range("H105").Value = 3944.2856
$ws.Range("K105").Value = 4101.6665
$ws.Range("H126").Value = 44809.832
$ws.Range("L126").Value = 44809.832
$ws.Range("J126").Value = 44809.832
$ws.Range("N126").Value = -54689.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L16").Value = 500
$ws.Range("H16").Value = 650
$ws.Range("J16").Value = 500
$ws.Range("N16").Value = -1074
$ws.Range("L31").Value = 2553.75
$ws.Range("K31").Value = 1254.6364
$ws.Range("I31").Value = 1254.6364
$ws.Range("J31").Value = 2553.75
$ws.Range("N31").Value = -3143.75
$ws.Range("M31").Value = -959.6364000000001
$ws.Range("H31").Value = 2024.4814
$ws.Range("N34").Value = -2957.75
$ws.Range("L34").Value = 2553.75
$ws.Range("M34").Value = -1052.6364
$ws.Range("I34").Value = 1254.6364
$ws.Range("J34").Value = 2553.75
$ws.Range("K34").Value = 1254.6364
$ws.Range("H34").Value = 2024.4814
$ws.Range("J41").Value = 14000
$ws.Range("H41").Value = 5600
$ws.Range("K41").Value = 2800
$ws.Range("N41").Value = -14856
$ws.Range("I41").Value = 2800
$ws.Range("M41").Value = -2372
$ws.Range("L41").Value = 14000
$ws.Range("L51").Value = 18643.334
$ws.Range("H51").Value = 18643.334
$ws.Range("J51").Value = 18643.334
$ws.Range("N51").Value = -20115.334
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -42290
$ws.Range("K59").Value = 2200
$ws.Range("M59").Value = -1055
$ws.Range("I59").Value = 2200
$ws.Range("H59").Value = 31600
$ws.Range("I60").Value = 6333.3335
$ws.Range("M60").Value = -5822.3335
$ws.Range("K60").Value = 6333.3335
$ws.Range("H60").Value = 6333.3335
$ws.Range("L61").Value = 18643.334
$ws.Range("N61").Value = -19339.334
$ws.Range("H61").Value = 18643.334
$ws.Range("J61").Value = 18643.334
$ws.Range("J68").Value = 20780
$ws.Range("L68").Value = 20780
$ws.Range("H68").Value = 20780
$ws.Range("N68").Value = -22278
$ws.Range("J71").Value = 20780
$ws.Range("N71").Value = -69828
$ws.Range("L71").Value = 62340
$ws.Range("H71").Value = 20780
$ws.Range("L74").Value = 25000
$ws.Range("H74").Value = 25000
$ws.Range("N74").Value = -26748
$ws.Range("J74").Value = 25000
$ws.Range("L77").Value = 75000
$ws.Range("H77").Value = 25000
$ws.Range("N77").Value = -83736
$ws.Range("J77").Value = 25000
$ws.Range("N113").Value = -4840
$ws.Range("H113").Value = 650
$ws.Range("J113").Value = 500
$ws.Range("L113").Value = 500
$ws.Range("L140").Value = 48238
$ws.Range("J140").Value = 48238
$ws.Range("H140").Value = 48238
$ws.Range("N140").Value = -58598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J68").Value = 1755.174
$ws.Range("L68").Value = 5265.522
$ws.Range("M68").Value = -1545.44
$ws.Range("K68").Value = 2356.44
$ws.Range("H68").Value = 1413.7324
$ws.Range("I68").Value = 785.48
$ws.Range("N68").Value = -6887.522
$ws.Range("M71").Value = -3013.32
$ws.Range("J71").Value = 1755.174
$ws.Range("N71").Value = -23908.566
$ws.Range("L71").Value = 15796.566
$ws.Range("I71").Value = 785.48
$ws.Range("K71").Value = 7069.32
$ws.Range("H71").Value = 1413.7324
$ws.Range("I134").Value = 7265.5713
$ws.Range("K134").Value = 21796.7139
$ws.Range("N134").Value = -34140
$ws.Range("M134").Value = -16726.7139
$ws.Range("L134").Value = 24000
$ws.Range("H134").Value = 7532.636
$ws.Range("J134").Value = 8000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J102").Value = 2719
$ws.Range("H102").Value = 1917.091
$ws.Range("L102").Value = 2719
$ws.Range("N102").Value = -5963
$ws.Range("K102").Value = 954.8
$ws.Range("M102").Value = 667.2
$ws.Range("I102").Value = 954.8
$ws.Range("J138").Value = 27500.818
$ws.Range("H138").Value = 27500.818
$ws.Range("L138").Value = 27500.818
$ws.Range("N138").Value = -37780.818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L7").Value = 1449.3334
$ws.Range("J7").Value = 1449.3334
$ws.Range("M7").Value = -984.25
$ws.Range("N7").Value = -1673.3334
$ws.Range("I7").Value = 1096.25
$ws.Range("H7").Value = 1247.5714
$ws.Range("K7").Value = 1096.25
$ws.Range("K22").Value = 250
$ws.Range("H22").Value = 1399.25
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 1526.9445
$ws.Range("M22").Value = 45
$ws.Range("N22").Value = -2116.9445
$ws.Range("L22").Value = 1526.9445
$ws.Range("K27").Value = 250
$ws.Range("N27").Value = -1740.9445
$ws.Range("M27").Value = -143
$ws.Range("H27").Value = 1399.25
$ws.Range("J27").Value = 1526.9445
$ws.Range("L27").Value = 1526.9445
$ws.Range("I27").Value = 250
$ws.Range("J55").Value = 394.66666
$ws.Range("M55").Value = 61.85714
$ws.Range("H55").Value = 222.08696
$ws.Range("K55").Value = 111.14286
$ws.Range("L55").Value = 394.66666
$ws.Range("I55").Value = 111.14286
$ws.Range("N55").Value = -740.66666
$ws.Range("K126").Value = 3288.75
$ws.Range("H126").Value = 1247.5714
$ws.Range("L126").Value = 4348.0002
$ws.Range("J126").Value = 1449.3334
$ws.Range("N126").Value = -9288.0002
$ws.Range("M126").Value = -818.75
$ws.Range("I126").Value = 1096.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 32558.385
$ws.Range("N123").Value = -42358.38499999999
$ws.Range("L123").Value = 32558.385
$ws.Range("J123").Value = 32558.385
$ws.Range("K126").Value = 1843.33338
$ws.Range("H126").Value = 1044.0834
$ws.Range("L126").Value = 6999
$ws.Range("J126").Value = 2333
$ws.Range("N126").Value = -11939
$ws.Range("M126").Value = 626.66662
$ws.Range("I126").Value = 614.44446
$ws.Range("I132").Value = 1194.579
$ws.Range("N132").Value = -14564.9
$ws.Range("L132").Value = 9504.900000000001
$ws.Range("M132").Value = -1053.737
$ws.Range("K132").Value = 3583.737
$ws.Range("J132").Value = 3168.3
$ws.Range("H132").Value = 1605.7709
$ws.Range("L133").Value = 59225
$ws.Range("H133").Value = 59225
$ws.Range("J133").Value = 59225
$ws.Range("N133").Value = -69345
$ws.Range("J138").Value = 46640
$ws.Range("H138").Value = 46640
$ws.Range("L138").Value = 46640
$ws.Range("N138").Value = -56920
